# Updates cryptos list data (prices and volume changes) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column (D2:D51) as Text so that numeric-looking
# values (e.g. "294.01") are preserved as strings, matching the source
# data which stores every Price/Volume cell as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.216.20"
$ws.Range("E2").Value = "  -8.85%  "

$ws.Range("D3").Value = "2.496.71"
$ws.Range("E3").Value = "  -4.38%  "

$ws.Range("D5").Value = "294.01"
$ws.Range("E5").Value = "  -4.42%  "

$ws.Range("D6").Value = "92.42"
$ws.Range("E6").Value = "  -7.29%  "

$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  -6.17%  "

$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.51%  "

$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  -6.71%  "

$ws.Range("D10").Value = "35.79"
$ws.Range("E10").Value = "  -9.32%  "

$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -6.10%  "

$ws.Range("D12").Value = "7.56"
$ws.Range("E12").Value = "  -6.92%  "

$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("D14").Value = "2.894.45"
$ws.Range("E14").Value = "  -3.72%  "

$ws.Range("D15").Value = "2.506.22"
$ws.Range("E15").Value = "  -4.14%  "

$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -6.95%  "

$ws.Range("D17").Value = "13.90"
$ws.Range("E17").Value = "  -6.97%  "

$ws.Range("D18").Value = "42.470.12"
$ws.Range("E18").Value = "  -8.64%  "

$ws.Range("D19").Value = "0.0₃0950"
$ws.Range("E19").Value = "  -6.04%  "

$ws.Range("D20").Value = "6.42"
$ws.Range("E20").Value = "  -4.53%  "

$ws.Range("D21").Value = "12.06"
$ws.Range("E21").Value = "  -6.55%  "

$ws.Range("D22").Value = "71.50"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "254.23"
$ws.Range("E23").Value = "  -6.62%  "

$ws.Range("D24").Value = "2.87"
$ws.Range("E24").Value = "  -5.85%  "

$ws.Range("D25").Value = "2.09"
$ws.Range("E25").Value = "  -3.55%  "

$ws.Range("D26").Value = "28.10"
$ws.Range("E26").Value = "  -3.48%  "

$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -2.82%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "9.86"
$ws.Range("E29").Value = "  -6.93%  "

$ws.Range("D30").Value = "36.36"
$ws.Range("E30").Value = "  -6.15%  "

$ws.Range("D31").Value = "5.88"
$ws.Range("E31").Value = "  -7.09%  "

$ws.Range("D32").Value = "3.42"
$ws.Range("E32").Value = "  -6.41%  "

$ws.Range("D33").Value = "150.46"
$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("D34").Value = "2.14"
$ws.Range("E34").Value = "  -4.32%  "

$ws.Range("E35").Value = "  -5.14%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0783"
$ws.Range("E36").Value = "  -6.55%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.112"
$ws.Range("E37").Value = "  -7.93%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "0.117"
$ws.Range("E38").Value = "  -4.75%  "

$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "23.54"
$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").Value = "16.45"
$ws.Range("E40").Value = "  +3.48%  "

$ws.Range("D41").Value = "3.40"
$ws.Range("E41").Value = "  -6.37%  "

$ws.Range("D42").Value = "0.0306"
$ws.Range("E42").Value = "  -7.60%  "

$ws.Range("D43").Value = "3.75"
$ws.Range("E43").Value = "  -7.79%  "

$ws.Range("D44").Value = "2.004.06"
$ws.Range("E44").Value = "  -5.15%  "

$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "1.64"
$ws.Range("E46").Value = "  +5.04%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "83.68"
$ws.Range("E47").Value = "  -10.44%  "

$ws.Range("D48").Value = "8.85"
$ws.Range("E48").Value = "  -7.48%  "

$ws.Range("D49").Value = "2.756.08"
$ws.Range("E49").Value = "  -3.82%  "

$ws.Range("D50").Value = "101.26"
$ws.Range("E50").Value = "  -7.04%  "

$ws.Range("D51").Value = "0.184"
$ws.Range("E51").Value = "  -8.05%  "

# Restore the default "Normal" style on the Price column so no stray
# cell-style index is left behind, keeping values stored as text.
$ws.Range("D2:D51").Style = "Normal"
